$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (columns A-T), per the revised NATMI
# output (3 clusters x 3 clusters, including self-pairs), as advised by
# Dr Hou.
$data = @(
    @("ECs","Efnb1","Erbb2","ECs",3,1,3.896501666666667,11.689505,0.401720501899026,0.401720501899026,2,0.6666666666666666,1.534538333333333,4.603615,0.1494637976135089,0.1494637976135089,5.979331173397222,53.813980560575,0.06004267179303326,0.06004267179303326),
    @("ECs","Efnb1","Erbb2","FAPs",3,1,3.896501666666667,11.689505,0.401720501899026,0.401720501899026,3,1,4.973328333333334,14.919985,0.4844014146353658,0.4844014146353658,19.37858213971389,174.407239257425,0.1945939794079173,0.1945939794079173),
    @("ECs","Efnb1","Erbb2","sCs",3,1,3.896501666666667,11.689505,0.401720501899026,0.401720501899026,3,1,3.75909,11.27727,0.3661347877511252,0.3661347877511252,14.64730045015,131.82570405135,0.1470838506980754,0.1470838506980754),
    @("FAPs","Efnb1","Erbb2","ECs",3,1,3.23724,9.71172,0.3337521163387849,0.3337521163387849,2,0.6666666666666666,1.534538333333333,4.603615,0.1494637976135089,0.1494637976135089,4.967668874199999,44.7090198678,0.04988385876954044,0.04988385876954044),
    @("FAPs","Efnb1","Erbb2","FAPs",3,1,3.23724,9.71172,0.3337521163387849,0.3337521163387849,3,1,4.973328333333334,14.919985,0.4844014146353658,0.4844014146353658,16.0998574138,144.8987167242,0.1616699972920546,0.1616699972920546),
    @("FAPs","Efnb1","Erbb2","sCs",3,1,3.23724,9.71172,0.3337521163387849,0.3337521163387849,3,1,3.75909,11.27727,0.3661347877511252,0.3661347877511252,12.1690765116,109.5216886044,0.1221982602771899,0.1221982602771899),
    @("sCs","Efnb1","Erbb2","ECs",3,1,2.565792333333333,7.697377,0.2645273817621892,0.2645273817621892,2,0.6666666666666666,1.534538333333333,4.603615,0.1494637976135089,0.1494637976135089,3.937306690872777,35.435760217855,0.03953726705093526,0.03953726705093526),
    @("sCs","Efnb1","Erbb2","FAPs",3,1,2.565792333333333,7.697377,0.2645273817621892,0.2645273817621892,3,1,4.973328333333334,14.919985,0.4844014146353658,0.4844014146353658,12.76052770881611,114.844749379345,0.1281374379353939,0.1281374379353939),
    @("sCs","Efnb1","Erbb2","sCs",3,1,2.565792333333333,7.697377,0.2645273817621892,0.2645273817621892,3,1,3.75909,11.27727,0.3661347877511252,0.3661347877511252,9.645044302310001,86.80539872079,0.09685267677585999,0.09685267677585999)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
